$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.413.19'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.008.76'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.95'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.65%  '
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0769'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.303.68'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.22'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.798'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.35%  '
$ws.Range('E16').Value = '  -5.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.005.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.369.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('E19').Value = '  -1.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0835'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.78%  '
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Value = '  +3.97%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.58'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('E29').Value = '  -4.94%  '
$ws.Range('E30').Value = '  -5.41%  '
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.60'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0643'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.51'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('E35').Value = '  -3.56%  '
$ws.Range('E36').Value = '  +0.52%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.35'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('E42').Value = '  -1.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0923'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.413.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('E45').Value = '  -5.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.07%  '
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.43%  '
$ws.Range('E49').Value = '  -6.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.195.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('E51').Value = '  -8.30%  '

Write-Host "Update complete"
